$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new sheets at the end of the workbook (after SignUpPage)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SignUpPage")
$wsProfile = $wb.Worksheets.Add($null, $ws3)
$wsProfile.Name = "ProfilePage"
$wsRFA = $wb.Worksheets.Add($null, $wsProfile)
$wsRFA.Name = "CreateRFA"

# ---------------------------------------------------------------------------
# 2. ErrorMessages sheet: update B8 text, add rows 10-12
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ErrorMessages")
$ws2.Range("B8").Value = "Company Name Already Exists. Please try With Different Name."
$ws2.Range("A10").Value = "CreateRFPError"
$ws2.Range("B10").Value = "Please Enter product name."
$ws2.Range("A11").Value = "CreateRFPSupplierError"
$ws2.Range("B11").Value = "Please select atleast one supplier."
$ws2.Range("A12").Value = " "

# ---------------------------------------------------------------------------
# 3. SignUpPage sheet: populate column C (mirrors column B, new hyperlink)
# ---------------------------------------------------------------------------
$ws3.Range("C2").Value = "Ajay"
$ws3.Range("C3").Value = "Singh"

$ws3.Range("B4").Copy()
$ws3.Range("C4").PasteSpecial(-4122) | Out-Null
$ws3.Range("C4").Value = "Sam Technology"

$ws3.Range("B5").Copy()
$ws3.Range("C5").PasteSpecial(-4122) | Out-Null
$ws3.Range("C5").Value = "Noida sector 62"

$ws3.Range("B6").Copy()
$ws3.Range("C6").PasteSpecial(-4122) | Out-Null
$ws3.Range("C6").Value = "Uttar Pradesh"

$ws3.Range("B7").Copy()
$ws3.Range("C7").PasteSpecial(-4122) | Out-Null
$ws3.Range("C7").Value = "Noida"

$ws3.Range("B8").Copy()
$ws3.Range("C8").PasteSpecial(-4122) | Out-Null
$ws3.Range("C8").Value = 201301

$ws3.Range("B9").Copy()
$ws3.Range("C9").PasteSpecial(-4122) | Out-Null
$ws3.Range("C9").Value = "01234AFD12"

$ws3.Range("B10").Copy()
$ws3.Range("C10").PasteSpecial(-4122) | Out-Null
$ws3.Range("C10").Value = "hprankit@gmail.com"
$ws3.Hyperlinks.Add($ws3.Range("C10"), "mailto:hprankit@gmail.com") | Out-Null
$ws3.Range("C10").Value = "hprankit@gmail.com"
$ws3.Range("B10").Copy()
$ws3.Range("C10").PasteSpecial(-4122) | Out-Null

$ws3.Range("B11").Copy()
$ws3.Range("C11").PasteSpecial(-4122) | Out-Null
$ws3.Range("C11").Value = 9568989975

$ws3.Range("B12").Copy()
$ws3.Range("C12").PasteSpecial(-4122) | Out-Null
$ws3.Range("C12").Value = "qwerty11"

# ---------------------------------------------------------------------------
# 4. ProfilePage sheet content
# ---------------------------------------------------------------------------
$wsProfile.Range("A1").Copy()
$wsProfile.Range("B1").PasteSpecial(-4122) | Out-Null
$wsProfile.Range("A1").Value = "objectID"
$wsProfile.Range("B1").Value = "MaxLengthValue"

$wsProfile.Range("A2").Value = "GstField"
$ws3.Range("B2").Copy()
$wsProfile.Range("B2").PasteSpecial(-4122) | Out-Null
$wsProfile.Range("B2").Value = 15

$wsProfile.Columns("B").ColumnWidth = 18.5546875

# ---------------------------------------------------------------------------
# 5. CreateRFA sheet content
# ---------------------------------------------------------------------------
$wsRFA.Range("A1").Value = "objectID"
$wsRFA.Range("B1").Value = "Category"
$wsRFA.Range("C1").Value = "SubCategory"
$wsRFA.Range("D1").Value = "SuppliersName"
$wsRFA.Range("A1:D1").Font.Bold = $true

$wsRFA.Range("A2").Value = "Mineral Water"
$wsRFA.Range("B2").Value = "Mineral Water"
$wsRFA.Range("C2").Value = "Drinking Water,Package,Packaging Machine,Treatment Equipments"
$wsRFA.Range("D2").Value = "SHIVAM ENTERPRISES"

$wsRFA.Range("A3").Value = "Travel & Hotels"
$wsRFA.Range("B3").Value = "Travel & Hotels "
$wsRFA.Range("C3").Value = "Guest Houses,OYO Rooms ,Travel & Lodging"
$wsRFA.Range("D3").Value = "ARCHER TOURS PVT LTD"

$wsRFA.Columns("A").ColumnWidth = 22.44140625
$wsRFA.Columns("B").ColumnWidth = 19.44140625
$wsRFA.Columns("C").ColumnWidth = 56.33203125
$wsRFA.Columns("D").ColumnWidth = 21.109375

$wsRFA.Range("C13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B9").Select() | Out-Null

$wsProfile.Activate()
$wsProfile.Range("A1:B1").Select() | Out-Null

$ws3.Activate()
$ws3.Range("B9").Select() | Out-Null
